$wb = $excel.ActiveWorkbook

# --- Update overall status text from "Ready for handoff" to "Handback transform failed" ---
# This status text is shared by the Overview sheet (row for the a2f8fb6d... file)
# and the per-language detail sheets (zh-cn, de-de) for the same file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"

# --- Fill in the "Error Detail" column (L) for the failed handback on each language sheet ---
$wsZhCn.Range("L3").Value = "Handback file name: e2ovjf2b.l1q is different with handoff file name: a2f8fb6d-4fb1-4486-b3ba-730fecd8a1df.6da92b75219e337fbf094efe5c8bc55dfbe25dfe.zh-cn."

$wsDeDe.Range("L3").Value = "Handback file name: e2ovjf2b.l1q is different with handoff file name: a2f8fb6d-4fb1-4486-b3ba-730fecd8a1df.6da92b75219e337fbf094efe5c8bc55dfbe25dfe.de-de."
